# The deck ships two themes:
#   theme1.xml -> "Office Theme" (reachable only via the Notes Master)
#   theme2.xml -> "Integral"     (the theme actually applied to the
#                                  slide master / whole presentation)
#
# The authored edit swaps the two themes' content, so the slides (which
# are driven by theme2.xml) switch from the "Integral" palette to the
# classic "Office Theme" palette (and vice-versa for the notes master's
# theme). The only theme that the PowerPoint object model lets us touch
# from here is the one bound to the slide master/presentation
# (theme2.xml) via ColorScheme / ThemeColorScheme, so we recolor it to
# the "Office Theme" swatch that should land there.
#
# Colors are expressed as VBA-style RGB() long values (0xBBGGRR) because
# that's what ThemeColor.RGB expects/returns.
function BGR([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Use the slide's ThemeColorScheme (rather than Master.ColorScheme) so the
# underlying clrScheme's metadata stays intact while we recolor it.
$tcs = $p.Slides.Item(1).ThemeColorScheme

# Target palette = the presentation's original "Office Theme" colors
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), in that slot order.
$tcs.Item(1).RGB  = BGR 0x00 0x00 0x00   # dk1      000000
$tcs.Item(2).RGB  = BGR 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Item(3).RGB  = BGR 0x44 0x54 0x6A   # dk2      44546A
$tcs.Item(4).RGB  = BGR 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Item(5).RGB  = BGR 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Item(6).RGB  = BGR 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Item(7).RGB  = BGR 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Item(8).RGB  = BGR 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Item(9).RGB  = BGR 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Item(10).RGB = BGR 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Item(11).RGB = BGR 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Item(12).RGB = BGR 0x95 0x4F 0x72   # folHlink 954F72
